# Rearrange the header labels so that multiple primary keys are handled:
# B1 -> Risk, C1 -> Curve, D1 -> Type (A1 stays TradeID)
# and move the active selection to B1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Risk"
$ws.Range("C1").Value = "Curve"
$ws.Range("D1").Value = "Type"

$ws.Range("B1").Select()
